$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the raw test-count values that drive the summary formulas.
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 23

# Move the active selection to A10 as recorded in the sheet view.
$ws.Range("A10").Select()

$excel.CalculateFullRebuild()
